# [PBL] Planning Sprint 4.
# Add a new "Sprint 4" row to the Sprints backlog sheet, and fill in the
# CapacityForecast/EffortForecast numbers for Sprint 3 (row 4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Finish filling in Sprint 3's forecast numbers (row 4).
$ws.Range("F4").Value = 8
$ws.Range("G4").Value = 3

# Add the new Sprint 4 row (row 5).
$ws.Range("A5").Value = "Sprint 4"

# Copy the date formatting from the row above so the new dates pick up
# the same number format/style as the existing Start/End date cells.
$ws.Range("B4:C4").Copy()
$ws.Range("B5:C5").PasteSpecial(-4122)

$ws.Range("B5").Value = 41757
$ws.Range("C5").Value = 41759
$ws.Range("D5").Value = 4
$ws.Range("E5").Value = 2

# Match the saved selection state recorded in the workbook.
$ws.Range("F5").Select()
